# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 16-41 of "Hoja1": (DocNumber, Name, Period, ValorMora)
$rows = @(
    @{ r = 16; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2008"; f = 35112 },
    @{ r = 17; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2009"; f = 35112 },
    @{ r = 18; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2010"; f = 35112 },
    @{ r = 19; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2011"; f = 35112 },
    @{ r = 20; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2012"; f = 35112 },
    @{ r = 21; c = "1047459128"; d = "ANDREA CAMILA GONZALEZ ARNEDO"; e = "2101"; f = 35112 },
    @{ r = 22; c = "1044927883"; d = "HERMOGENES ZAPATEIRO PAJARO";   e = "2101"; f = 35112 },
    @{ r = 23; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2101"; f = 35112 },
    @{ r = 24; c = "1047459128"; d = "ANDREA CAMILA GONZALEZ ARNEDO"; e = "2102"; f = 35112 },
    @{ r = 25; c = "1044927883"; d = "HERMOGENES ZAPATEIRO PAJARO";   e = "2102"; f = 35112 },
    @{ r = 26; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2102"; f = 35112 },
    @{ r = 27; c = "1047459128"; d = "ANDREA CAMILA GONZALEZ ARNEDO"; e = "2103"; f = 35112 },
    @{ r = 28; c = "1044927883"; d = "HERMOGENES ZAPATEIRO PAJARO";   e = "2103"; f = 35112 },
    @{ r = 29; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2103"; f = 35112 },
    @{ r = 30; c = "1047459128"; d = "ANDREA CAMILA GONZALEZ ARNEDO"; e = "2104"; f = 35112 },
    @{ r = 31; c = "1044927883"; d = "HERMOGENES ZAPATEIRO PAJARO";   e = "2104"; f = 35112 },
    @{ r = 32; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2104"; f = 35112 },
    @{ r = 33; c = "1047459128"; d = "ANDREA CAMILA GONZALEZ ARNEDO"; e = "2105"; f = 35112 },
    @{ r = 34; c = "1044927883"; d = "HERMOGENES ZAPATEIRO PAJARO";   e = "2105"; f = 35112 },
    @{ r = 35; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2105"; f = 35112 },
    @{ r = 36; c = "1047459128"; d = "ANDREA CAMILA GONZALEZ ARNEDO"; e = "2106"; f = 35112 },
    @{ r = 37; c = "1044927883"; d = "HERMOGENES ZAPATEIRO PAJARO";   e = "2106"; f = 35112 },
    @{ r = 38; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2106"; f = 35112 },
    @{ r = 39; c = "1047459128"; d = "ANDREA CAMILA GONZALEZ ARNEDO"; e = "2107"; f = 29260 },
    @{ r = 40; c = "1044927883"; d = "HERMOGENES ZAPATEIRO PAJARO";   e = "2107"; f = 29260 },
    @{ r = 41; c = "73572193";   d = "MIGUEL DIONISIO LEON DALMAU";   e = "2107"; f = 29260 }
)

foreach ($row in $rows) {
    $ws.Range("C" + $row.r).Value = $row.c
    $ws.Range("D" + $row.r).Value = $row.d
    $ws.Range("E" + $row.r).Value = $row.e
    $ws.Range("F" + $row.r).Value = $row.f
}
